$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "Price" column values keep their exact textual representation
# (e.g. trailing zeros, thousands-style dots) by forcing Text format before
# assigning, then update both the Price (D) and Volume(1h) (E) columns.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.923.16"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.812.23"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9980"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.45"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9988"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4994"
$ws.Range("E7").Value = "  -4.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3893"
$ws.Range("E8").Value = "  +1.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09544"
$ws.Range("E9").Value = "  +20.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.101"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.34"
$ws.Range("E11").Value = "  -2.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.403"
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9978"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.52"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.813.41"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.253"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001127"
$ws.Range("E17").Value = "  +3.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.16"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06584"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9997"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.19"
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.952"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.955.15"
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.243"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.76"
$ws.Range("E26").Value = "  +1.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.12"
$ws.Range("E27").Value = "  -1.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.019.10"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.410"
$ws.Range("E29").Value = "  +3.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.44"
$ws.Range("E30").Value = "  +4.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1073"
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.058"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.636"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.618"
$ws.Range("E34").Value = "  -1.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06815"
$ws.Range("E35").Value = "  -7.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.968"
$ws.Range("E36").Value = "  +2.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02313"
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2155"
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.41"
$ws.Range("E39").Value = "  -6.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.957"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6249"
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.145"
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.07"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5900"
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.286"
$ws.Range("E46").Value = "  -6.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.681"
$ws.Range("E47").Value = "  -2.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.27"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.959"
$ws.Range("E49").Value = "  +2.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.179"
$ws.Range("E50").Value = "  -4.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06763"
$ws.Range("E51").Value = "  +0.02%  "
